# "Search functions in the DAOs" — the search() methods have now been
# implemented, so every remaining task on the to-do list is done.
# Flip the few "Statut" cells that weren't already "terminé" to "terminé".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("D15").Value = "terminé"
$ws.Range("D21").Value = "terminé"
$ws.Range("D27").Value = "terminé"

# The Total row's Statut cell was styled like a header (bold). Make it
# match the rest of the (now all "terminé") Statut column: regular weight.
$ws.Range("D59").Value = "terminé"
$ws.Range("D59").Font.Bold = $false

# Leave the cursor where the author ended up after the edit.
$ws.Range("E48").Select()
